# Apply the fixes described in the commit "Fixed a few mistakes in datasets
# and rerun the whole code" to Reduced_with-real-values-community.xlsx.
#
# Summary of the edit:
#   1. Rename the "C_conversion_ratio" column header to "C_ratio".
#   2. Correct a handful of standalone data-entry mistakes
#      (AM25, W27, W28, AM48, AM132).
#   3. Remove the erroneous experiment that used to live in row 148
#      (Kein 2013 EPSL, experiment #146) - this shifts every following
#      row up by one and the recomputed ("rerun") derived columns for
#      those rows take on the values that used to belong to the next
#      row down. The running index in column A is NOT data that moves
#      with the row -- it is a recomputed sequential row counter, so it
#      is rebuilt from scratch after the deletion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text fix -------------------------------------------------
$ws.Range("AL1").Value = "C_ratio"

# --- 2. Standalone value corrections ------------------------------------
$ws.Range("AM25").Value = 1
$ws.Range("W27").Value = 30
$ws.Range("W28").Value = 30
$ws.Range("AM48").Value = 1
$ws.Range("AM132").Value = 1

# --- 3. Drop the bad experiment row and shift everything below it up ---
$ws.Rows("148").Delete()

# Column A is a plain sequential row counter (row number - 2); rebuild it
# for every data row now that row 148 is gone.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 158) { $lastRow = 158 }
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
